$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column C header / values ---
$ws.Range("C1").Value = "Folder"
$ws.Range("C2").Value = "Full_Install"
$ws.Range("C3").Value = "Custom_Release"
$ws.Range("C4").Value = "Dictionary_Metadata_EditClar"
$ws.Range("C5").Value = "Dictionary_Clinical_Notes"
$ws.Range("C6").Value = "Dictionary_Only"

# --- Column C width ---
$ws.Range("C1").ColumnWidth = 92

# --- Set active cell / selection ---
$ws.Range("C10").Select()
